$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.599.75"
$ws.Range("E2").Value = "  +4.35%  "

$ws.Range("D3").Value = "'4.043.76"
$ws.Range("E3").Value = "  +3.17%  "

$ws.Range("D5").Value = "'517.47"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").Value = "'147.08"
$ws.Range("E6").Value = "  +1.56%  "

$ws.Range("D7").Value = "'0.721"
$ws.Range("E7").Value = "  +17.38%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.763"
$ws.Range("E9").Value = "  +5.40%  "

$ws.Range("D10").Value = "'0.175"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").Value = "'0.0000324"
$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("D12").Value = "'46.54"
$ws.Range("E12").Value = "  +9.97%  "

$ws.Range("D13").Value = "'10.85"
$ws.Range("E13").Value = "  +5.38%  "

$ws.Range("D14").Value = "'4.676.37"
$ws.Range("E14").Value = "  +2.91%  "

$ws.Range("D15").Value = "'4.040.47"
$ws.Range("E15").Value = "  +3.10%  "

$ws.Range("D16").Value = "'21.16"
$ws.Range("E16").Value = "  +6.79%  "

$ws.Range("D17").Value = "'14.13"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'72.282.99"
$ws.Range("E20").Value = "  +3.83%  "

$ws.Range("D21").Value = "'442.52"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("D22").Value = "'104.57"
$ws.Range("E22").Value = "  +17.74%  "

$ws.Range("D23").Value = "'3.55"
$ws.Range("E23").Value = "  +5.18%  "

$ws.Range("D24").Value = "'14.67"
$ws.Range("E24").Value = "  +2.96%  "

$ws.Range("D25").Value = "'3.97"
$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("D26").Value = "'11.49"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").Value = "'11.13"
$ws.Range("E27").Value = "  +4.04%  "

$ws.Range("D28").Value = "'38.31"
$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("E29").Value = "  +2.03%  "

$ws.Range("D30").Value = "'3.10"
$ws.Range("E30").Value = "  +9.96%  "

$ws.Range("D31").Value = "'13.72"
$ws.Range("E31").Value = "  +3.90%  "

$ws.Range("D32").Value = "'678.18"
$ws.Range("E32").Value = "  -2.34%  "

$ws.Range("D34").Value = "'6.78"
$ws.Range("E34").Value = "  +13.12%  "

$ws.Range("D35").Value = "'67.39"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").Value = "'41.45"
$ws.Range("E36").Value = "  +3.62%  "

$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "'0.0₃0862"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.429"
$ws.Range("E38").Value = "  -3.05%  "

$ws.Range("D39").Value = "'3.52"
$ws.Range("E39").Value = "  +14.24%  "

$ws.Range("D40").Value = "'0.150"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").Value = "'0.0491"
$ws.Range("E43").Value = "  +1.78%  "

$ws.Range("D44").Value = "'3.18"
$ws.Range("E44").Value = "  +3.06%  "

$ws.Range("E45").Value = "  +11.21%  "

$ws.Range("E46").Value = "  -2.16%  "

$ws.Range("E47").Value = "  +4.76%  "

$ws.Range("E48").Value = "  +1.06%  "

$ws.Range("D49").Value = "'9.17"
$ws.Range("E49").Value = "  +8.40%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'2.11"
$ws.Range("E50").Value = "  +2.15%  "

$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.31"
$ws.Range("E51").Value = "  +0.14%  "

